$wb = $excel.ActiveWorkbook

# --- Rename the "Include from LOINC" sheet to "Include #0" ---
$wsInclude = $wb.Worksheets.Item("Include from LOINC")
$wsInclude.Name = "Include #0"

# --- Update the Metadata sheet ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# Update the Date property value (row 8, column B)
$wsMeta.Range("B8").Value = "2024-09-17T19:55:11+00:00"

# Insert a new "Jurisdiction" property row after "Contact" (row 10),
# pushing Description/Purpose/Copyright/Immutable down by one row.
$wsMeta.Rows.Item(11).Insert()

# Match the formatting of the surrounding property rows.
$wsMeta.Range("A10:B10").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new row's values ("Jurisdiction" with an empty value).
$wsMeta.Range("A11").Value = "Jurisdiction"
